$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SPL", 4, 5, 4.262249820071284, 3.333893823378794, 0.00974046742299407, 0.624031007751938, 57.14285714285714),
    @("PEP", 2, 2, 0.2695588772136407, 0.2035406185491022, 0.009050882111654309, 0.5465116279069767, 66.40926640926641),
    @("KRU", 0, 5, 2.618842383902459, 2.001682927551309, 0.01239655736930432, 0.6162790697674418, 50.57915057915058),
    @("NEU", 8, 7, 3.581292169103917, 2.662951816310947, 0.008951901977750894, 0.5542635658914729, 64.09266409266409),
    @("ERB", 2, 6, 0.2852847340679718, 0.2023870455546976, 0.01512462150663506, 0.5155038759689923, 44.01544401544402),
    @("ATD", 1, 2, 0.0483899941386483, 0.03401734566914823, 0.009197572952402472, 0.3798449612403101, 66.02316602316603),
    @("DAT", 3, 1, 1.522083412802567, 1.021068681763854, 0.02907598781905718, 0.3488372093023256, 25.86872586872587),
    @("ZMT", 2, 3, 0.02576193131210118, 0.0147605819817843, 0.01730472196972008, 0.1937984496124031, 60.61776061776062),
    @("KTY", 1, 1, 4.980710419265365, 3.693641805455617, 0.01117038338088966, 0.4651162790697674, 57.14285714285714),
    @("ABE", 10, 2, 0.3728079688391112, 0.2667566143420379, 0.0136428731980587, 0.4224806201550387, 49.03474903474903),
    @("MRB", 0, 8, 0.01841025824737029, 0.01405469968690697, 0.013908950877897, 0.4573643410852713, 46.71814671814672),
    @("UNI", 10, 8, 0.0693466576127281, 0.05262296550380359, 0.008515165414469669, 0.5503875968992248, 66.40926640926641),
    @("WIG", 3, 2, 548.6331449524837, 439.626471200701, 0.007561880183372149, 0.7015503875968992, 73.35907335907336),
    @("WIG20", 2, 2, 15.96880173581907, 12.38793411847845, 0.005499783806024364, 0.6511627906976745, 85.71428571428571),
    @("mWIG40", 3, 4, 21.18662715589623, 16.47272245773662, 0.004148136970349986, 0.6744186046511628, 93.43629343629344),
    @("sWIG80", 7, 9, 221.5425959676884, 187.5055945082987, 0.01599344015247646, 0.7209302325581395, 32.81853281853282),
)


$startRow = 20

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i

    # Copy the style of A19 (bold, bordered, centered => style index 1)
    # onto the new column-A cell so formatting matches the existing rows.
    $ws.Range("A19").Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null

    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Cells.Item($r, 2).Value = "ARIMAX"
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $ws.Cells.Item($r, 9).Value = $row[6]
    $ws.Cells.Item($r, 10).Value = $row[7]
}

$excel.CutCopyMode = 0
